$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "coa" column with header and values
$ws.Range("E1").Value = "coa"
$ws.Range("E2").Value = "102.18.000"
$ws.Range("E3").Value = "102.18.000"

# Update selection to match the diff (E6 selected, nothing in it)
$ws.Range("E6").Select()
